# Update "DSM Scheduled Flights vs actual" worksheet:
#  - correct C1191 (85 -> 86) and restyle A1191:C1191 to the "text date / compact
#    number" look used by the newly appended rows
#  - append 20 new daily rows (2023-07-12 .. 2023-07-31) with scheduled /
#    tracked flight counts and the Percent formula carried down column D
#  - move the window scroll/selection to reflect where the user left off

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row r -> (DateText, Scheduled, Tracked)
$data = @(
    @(1191, "2023-07-11", 87, 86),
    @(1192, "2023-07-12", 74, 61),
    @(1193, "2023-07-13", 84, 83),
    @(1194, "2023-07-14", 80, 73),
    @(1195, "2023-07-15", 75, 71),
    @(1196, "2023-07-16", 68, 66),
    @(1197, "2023-07-17", 71, 68),
    @(1198, "2023-07-18", 70, 70),
    @(1199, "2023-07-19", 90, 86),
    @(1200, "2023-07-20", 91, 85),
    @(1201, "2023-07-21", 87, 86),
    @(1202, "2023-07-22", 68, 66),
    @(1203, "2023-07-23", 65, 64),
    @(1204, "2023-07-24", 76, 71),
    @(1205, "2023-07-25", 74, 71),
    @(1206, "2023-07-26", 74, 71),
    @(1207, "2023-07-27", 85, 77),
    @(1208, "2023-07-28", 81, 77),
    @(1209, "2023-07-29", 78, 74),
    @(1210, "2023-07-30", 64, 60),
    @(1211, "2023-07-31", 72, 68)
)

$firstRow = 1191
$lastRow = 1211

# --- column A: switch from date-serial cells to plain text "yyyy-mm-dd" labels,
#     bold / centered / wrapped (matches the rest of the appended block) -------
$colA = $ws.Range("A$firstRow`:A$lastRow")
$colA.NumberFormat = "@"
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4108
$colA.WrapText = $true
$colA.Font.Bold = $true
$colA.Font.Size = 10

# --- columns B & C: compact integer style used for the rest of the table -----
$colBC = $ws.Range("B$firstRow`:C$lastRow")
$colBC.NumberFormat = "0"
$colBC.WrapText = $true
$colBC.Font.Bold = $false
$colBC.Font.Size = 10

# --- column D: percent formula, formatted like the rest of the column --------
$colD = $ws.Range("D$firstRow`:D$lastRow")
$colD.NumberFormat = "0.0%"

foreach ($row in $data) {
    $r = $row[0]
    $dateText = $row[1]
    $scheduled = $row[2]
    $tracked = $row[3]

    $ws.Range("A$r").Value = $dateText
    $ws.Range("B$r").Value = $scheduled
    $ws.Range("C$r").Value = $tracked
}

$ws.Range("D$firstRow`:D$lastRow").Formula = "=C$firstRow/B$firstRow"

# --- reposition the view like the author left it ------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 635
$win.ScrollColumn = 1
$ws.Range("D1191:D1211").Select()

Write-Output "Rows $firstRow to $lastRow updated."
